$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edits 1, 2 & 4: collapse the split "<id>" / "<xxx>" / "</id>" runs into a
# single run (the resulting run keeps the formatting of the first/last run,
# i.e. Courier New / color 7f6000 / size 9, which is exactly what a
# same-paragraph Find & Replace over the combined text produces).
# ---------------------------------------------------------------------------
$idReplacements = @(
    "<id>p133v_3</id>",
    "<id>p134r_1</id>",
    "<id>p134r_2</id>"
)

foreach ($txt in $idReplacements) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, $txt, 2)
    if (-not $ok) {
        Write-Output ("WARNING: could not find/replace " + $txt)
    }
}

# ---------------------------------------------------------------------------
# Edit 3: "<head>To mend a pierced mold</head>" -> "<head>For mending a
# pierced mold</head>", reshaping the run boundaries so that:
#   "For"      -> no explicit color (inherits automatic)
#   " "        -> color 000000
#   "mending"  -> no explicit color (inherits automatic)
#   " a pierced mold" -> unchanged (color 000000)
# ---------------------------------------------------------------------------
$headRng = $d.Content
$found = $headRng.Find.Execute("<head>To mend a pierced mold</head>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Output "WARNING: could not locate the <head>To mend...</head> paragraph"
} else {
    $base = $headRng.Start

    # "<head>" is 6 characters; "To mend" begins right after it.
    # "To " occupies [base+6, base+9); "mend" occupies [base+9, base+13).

    # Step 1: turn "mend" (a run with no explicit color) into "For" - since
    # the whole replaced range lies inside that single run, the new text
    # keeps that run's (colorless) formatting.  "mend" occupied
    # [base+9, base+13); the replacement "For" now occupies [base+9, base+12).
    $mendRange = $d.Range($base + 9, $base + 13)
    if ($mendRange.Text -ne "mend") {
        Write-Output ("WARNING: unexpected text at mend offset: [" + $mendRange.Text + "]")
    }
    $mendRange.Text = "For"

    # Step 2: append "mending" right after the newly created "For" run
    # (which now ends at base+12); that boundary's preceding run is the
    # colorless "For" run, so the inserted text inherits the same
    # (colorless) formatting.
    $afterFor = $d.Range($base + 12, $base + 12)
    $afterFor.InsertBefore("mending")

    # Step 3: remove the original "To " (still at the front, colored black).
    $toRange = $d.Range($base + 6, $base + 9)
    if ($toRange.Text -ne "To ") {
        Write-Output ("WARNING: unexpected text at To offset: [" + $toRange.Text + "]")
    }
    $toRange.Delete()

    # Step 4: insert the separating space between "For" (now [base+6,base+9))
    # and "mending" (now [base+9,base+16)).
    $spacePos = $d.Range($base + 9, $base + 9)
    $spacePos.InsertBefore(" ")

    # Step 5: give that new space explicit black coloring (000000) to match
    # the surrounding body-text runs.
    $spaceRange = $d.Range($base + 9, $base + 10)
    if ($spaceRange.Text -ne " ") {
        Write-Output ("WARNING: unexpected text at space offset: [" + $spaceRange.Text + "]")
    }
    $spaceRange.Font.Color = 0

    $verify = $d.Range($base, $base + 39)
    Write-Output ("Head paragraph now: [" + $verify.Text + "]")
}

Write-Output "done"
